$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the style used for the date column (A) is available for the new row 19
# by copying the format from an existing date cell before writing the new row's values.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A1").Select()

$data = @(
    @(39400, 2007, 4.930115226412335, 2008, 1.402909115503936),
    @(39765, 2008, 1.457587285166628, 2009, 2.807231216534278),
    @(40130, 2009, -0.9140166223623569, 2010, 1.821983295885121),
    @(40494, 2010, 2.585942866987878, 2011, 2.722861752007866),
    @(40862, 2011, 4.253963781362402, 2012, 2.878414118480799),
    @(41228, 2012, 1.752870900283909, 2013, 3.144721336271927),
    @(41592, 2013, -1.479696720105139, 2014, 2.238623952069552),
    @(41957, 2014, 3.900127535411246, 2015, -0.6155071485167585),
    @(42321, 2015, 0.03947433952959933, 2016, 1.182212550358064),
    @(42689, 2016, 2.192778679161944, 2017, -0.5835597102573198),
    @(43053, 2017, 3.40836448860673, 2018, 2.34299484087257),
    @(43418, 2018, 2.799070570134488, 2019, 3.056075254339996),
    @(43783, 2019, 4.195393191694419, 2020, 2.652948310315506),
    @(44159, 2020, 1.666553973046048, 2021, -0.9999522486825452),
    @(44525, 2021, 1.879266440112803, 2022, -0.1343977949472275),
    @(44890, 2022, -2.620683231370946, 2023, -2.83913279674276),
    @(45254, 2023, -3.036556262700274, 2024, -1.803491225663911),
    @(45618, 2024, -2.953443685011514, 2025, -2.75492543068685)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
